$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo.Range("D7").Value = 475.2
$wsGrupo.Range("M7").Value = 978.48
$wsGrupo.Range("K9").Value = 1218.24
$wsGrupo.Range("K10").Value = 609.12
$wsGrupo.Range("K15").Value = 812.16
$wsGrupo.Range("M53").Value = 191.81
$wsGrupo.Range("K56").Value = "11 de 54"
$wsGrupo.Range("M56").Value = "16 de 54"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual.Range("F7").Value = 1453.68
$wsMensual.Range("F9").Value = 3172.12
$wsMensual.Range("F10").Value = 1726.02
$wsMensual.Range("F15").Value = 812.16
$wsMensual.Range("F55").Value = 264.06
$wsMensual.Range("F56").Value = 264.06
$wsMensual.Range("F60").Value = 95230.84999999999

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumpl.Range("D3").Value = 5598.72
$wsCumpl.Range("E3").Value = 12070.4270988183
$wsCumpl.Range("F3").Value = 0.3168641909362133

$wsCumpl.Range("D10").Value = 11998.17
$wsCumpl.Range("E10").Value = -8117.09016465608
$wsCumpl.Range("F10").Value = 3.091451479749524

$wsCumpl.Range("D12").Value = 53275.31
$wsCumpl.Range("E12").Value = -612.1899999999951
$wsCumpl.Range("F12").Value = 1.011624643583593

$wsCumpl.Range("D14").Value = 89868.22
$wsCumpl.Range("E14").Value = 9148.286611906138
$wsCumpl.Range("F14").Value = 0.9076084692851999
